# "New crime data collected" - weekly CompStat data refresh.
# Updates the report header (volume / week-covering date range) and the
# Week-to-Date, 28-Day, Year-to-Date and 2-Year crime figures table
# (rows 14-30, columns C:N) on the CompStat_1 sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Report header: volume number and the week-covering date range ---
$ws.Range("A8").Value2 = "Volume 30   Number  14"
$ws.Range("C9").Value2 = "Report Covering the Week  4/3/2023  Through  4/9/2023"

# --- Weekly crime statistics table (rows 14-30) ---

# Row 14 - Murder
$ws.Range("D14").Value2 = 8
$ws.Range("E14").Value2 = -12.5
$ws.Range("F14").Value2 = 30
$ws.Range("H14").Value2 = 0
$ws.Range("I14").Value2 = 100
$ws.Range("J14").Value2 = 111
$ws.Range("K14").Value2 = -9.909909909909
$ws.Range("L14").Value2 = -13.043478260869
$ws.Range("M14").Value2 = -21.875
$ws.Range("N14").Value2 = -80.916030534351

# Row 15 - Rape
$ws.Range("C15").Value2 = 33
$ws.Range("D15").Value2 = 31
$ws.Range("E15").Value2 = 6.451612903225
$ws.Range("F15").Value2 = 119
$ws.Range("G15").Value2 = 107
$ws.Range("H15").Value2 = 11.214953271028
$ws.Range("I15").Value2 = 408
$ws.Range("J15").Value2 = 441
$ws.Range("K15").Value2 = -7.482993197278
$ws.Range("L15").Value2 = 8.222811671087
$ws.Range("M15").Value2 = 26.708074534161
$ws.Range("N15").Value2 = -49.877149877149

# Row 16 - Robbery
$ws.Range("C16").Value2 = 255
$ws.Range("D16").Value2 = 339
$ws.Range("E16").Value2 = -24.778761061946
$ws.Range("F16").Value2 = 1138
$ws.Range("G16").Value2 = 1215
$ws.Range("H16").Value2 = -6.337448559670
$ws.Range("I16").Value2 = 4102
$ws.Range("J16").Value2 = 4251
$ws.Range("K16").Value2 = -3.505057633498
$ws.Range("L16").Value2 = 43.929824561403
$ws.Range("M16").Value2 = -13.459915611814
$ws.Range("N16").Value2 = -82.026902685887

# Row 17 - Fel. Assault
$ws.Range("C17").Value2 = 497
$ws.Range("D17").Value2 = 508
$ws.Range("E17").Value2 = -2.165354330708
$ws.Range("F17").Value2 = 1979
$ws.Range("G17").Value2 = 1946
$ws.Range("H17").Value2 = 1.695786228160
$ws.Range("I17").Value2 = 6742
$ws.Range("J17").Value2 = 6216
$ws.Range("K17").Value2 = 8.462033462033
$ws.Range("L17").Value2 = 30.861801242236
$ws.Range("M17").Value2 = 61.523718255869
$ws.Range("N17").Value2 = -31.183015208737

# Row 18 - Burglary
$ws.Range("C18").Value2 = 239
$ws.Range("D18").Value2 = 331
$ws.Range("E18").Value2 = -27.794561933534
$ws.Range("F18").Value2 = 1047
$ws.Range("G18").Value2 = 1246
$ws.Range("H18").Value2 = -15.971107544141
$ws.Range("I18").Value2 = 3878
$ws.Range("J18").Value2 = 4206
$ws.Range("K18").Value2 = -7.798383262006
$ws.Range("L18").Value2 = 22.993973993022
$ws.Range("M18").Value2 = -18.289085545722
$ws.Range("N18").Value2 = -85.389194484213

# Row 19 - Gr. Larceny
$ws.Range("C19").Value2 = 914
$ws.Range("D19").Value2 = 904
$ws.Range("E19").Value2 = 1.106194690265
$ws.Range("F19").Value2 = 3715
$ws.Range("G19").Value2 = 3664
$ws.Range("H19").Value2 = 1.391921397379
$ws.Range("I19").Value2 = 12839
$ws.Range("J19").Value2 = 13161
$ws.Range("K19").Value2 = -2.446622597067
$ws.Range("L19").Value2 = 53.815742182820
$ws.Range("M19").Value2 = 38.650107991360
$ws.Range("N19").Value2 = -38.911357472522

# Row 20 - G.L.A.
$ws.Range("C20").Value2 = 287
$ws.Range("D20").Value2 = 235
$ws.Range("E20").Value2 = 22.127659574468
$ws.Range("F20").Value2 = 1074
$ws.Range("G20").Value2 = 922
$ws.Range("H20").Value2 = 16.485900216919
$ws.Range("I20").Value2 = 3858
$ws.Range("J20").Value2 = 3566
$ws.Range("K20").Value2 = 8.188446438586
$ws.Range("L20").Value2 = 91.179385530227
$ws.Range("M20").Value2 = 51.531814611154
$ws.Range("N20").Value2 = -87.329633157082

# Row 21 - TOTAL
$ws.Range("C21").Value2 = 2232
$ws.Range("D21").Value2 = 2356
$ws.Range("E21").Value2 = -5.263157894736
$ws.Range("F21").Value2 = 9102
$ws.Range("G21").Value2 = 9130
$ws.Range("H21").Value2 = -0.306681270536
$ws.Range("I21").Value2 = 31927
$ws.Range("J21").Value2 = 31952
$ws.Range("K21").Value2 = -0.078242363545
$ws.Range("L21").Value2 = 45.043612574959
$ws.Range("M21").Value2 = 23.194165766322
$ws.Range("N21").Value2 = -71.485093689155

# Row 22 - Transit
$ws.Range("C22").Value2 = 39
$ws.Range("D22").Value2 = 34
$ws.Range("E22").Value2 = 14.705882352941
$ws.Range("F22").Value2 = 189
$ws.Range("G22").Value2 = 151
$ws.Range("H22").Value2 = 25.165562913907
$ws.Range("I22").Value2 = 578
$ws.Range("J22").Value2 = 621
$ws.Range("K22").Value2 = -6.924315619967
$ws.Range("L22").Value2 = 57.923497267759
$ws.Range("M22").Value2 = 3.956834532374

# Row 23 - Housing
$ws.Range("C23").Value2 = 109
$ws.Range("D23").Value2 = 105
$ws.Range("E23").Value2 = 3.809523809523
$ws.Range("F23").Value2 = 482
$ws.Range("H23").Value2 = 14.761904761904
$ws.Range("I23").Value2 = 1615
$ws.Range("J23").Value2 = 1490
$ws.Range("K23").Value2 = 8.389261744966
$ws.Range("L23").Value2 = 22.163388804841
$ws.Range("M23").Value2 = 65.641025641025

# Row 24 - Petit Larceny
$ws.Range("C24").Value2 = 1900
$ws.Range("D24").Value2 = 2233
$ws.Range("E24").Value2 = -14.912673533363
$ws.Range("F24").Value2 = 7793
$ws.Range("G24").Value2 = 8658
$ws.Range("H24").Value2 = -9.990759990759
$ws.Range("I24").Value2 = 28440
$ws.Range("J24").Value2 = 28310
$ws.Range("K24").Value2 = 0.459201695513
$ws.Range("L24").Value2 = 40.959555908009
$ws.Range("M24").Value2 = 44.732824427480

# Row 25 - Misd. Assault
$ws.Range("C25").Value2 = 832
$ws.Range("D25").Value2 = 801
$ws.Range("E25").Value2 = 3.870162297128
$ws.Range("F25").Value2 = 3343
$ws.Range("G25").Value2 = 3192
$ws.Range("H25").Value2 = 4.730576441102
$ws.Range("I25").Value2 = 10995
$ws.Range("J25").Value2 = 10377
$ws.Range("K25").Value2 = 5.955478461983
$ws.Range("L25").Value2 = 37.351655215490
$ws.Range("M25").Value2 = -3.535708018950

# Row 26 - UCR Rape*
$ws.Range("C26").Value2 = 56
$ws.Range("D26").Value2 = 51
$ws.Range("E26").Value2 = 9.803921568627
$ws.Range("F26").Value2 = 203
$ws.Range("G26").Value2 = 177
$ws.Range("H26").Value2 = 14.689265536723
$ws.Range("I26").Value2 = 658
$ws.Range("J26").Value2 = 708
$ws.Range("K26").Value2 = -7.062146892655
$ws.Range("L26").Value2 = 4.777070063694

# Row 27 - Other Sex Crimes
$ws.Range("C27").Value2 = 97
$ws.Range("D27").Value2 = 93
$ws.Range("E27").Value2 = 4.301075268817
$ws.Range("F27").Value2 = 407
$ws.Range("G27").Value2 = 409
$ws.Range("H27").Value2 = -0.488997555012
$ws.Range("I27").Value2 = 1327
$ws.Range("J27").Value2 = 1252
$ws.Range("K27").Value2 = 5.990415335463
$ws.Range("L27").Value2 = 22.530009233610

# Row 28 - Shooting Vic.
$ws.Range("C28").Value2 = 16
$ws.Range("D28").Value2 = 30
$ws.Range("E28").Value2 = -46.666666666666
$ws.Range("F28").Value2 = 89
$ws.Range("G28").Value2 = 125
$ws.Range("H28").Value2 = -28.8
$ws.Range("I28").Value2 = 292
$ws.Range("J28").Value2 = 359
$ws.Range("K28").Value2 = -18.662952646239
$ws.Range("L28").Value2 = -11.515151515151
$ws.Range("M28").Value2 = -24.155844155844
$ws.Range("N28").Value2 = -80.136054421768

# Row 29 - Shooting Inc.
$ws.Range("C29").Value2 = 15
$ws.Range("D29").Value2 = 24
$ws.Range("E29").Value2 = -37.5
$ws.Range("F29").Value2 = 78
$ws.Range("G29").Value2 = 107
$ws.Range("H29").Value2 = -27.102803738317
$ws.Range("I29").Value2 = 244
$ws.Range("J29").Value2 = 317
$ws.Range("K29").Value2 = -23.028391167192
$ws.Range("L29").Value2 = -17.288135593220
$ws.Range("M29").Value2 = -23.028391167192
$ws.Range("N29").Value2 = -81.858736059479

# Row 30 - Hate Crimes
$ws.Range("C30").Value2 = 4
$ws.Range("E30").Value2 = -50
$ws.Range("F30").Value2 = 42
$ws.Range("G30").Value2 = 39
$ws.Range("H30").Value2 = 7.692307692307
$ws.Range("I30").Value2 = 117
$ws.Range("J30").Value2 = 200
$ws.Range("K30").Value2 = -41.5
$ws.Range("L30").Value2 = 17
